# simplify steel description (remove RME)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared-string text in B2: drop the "/RME" segment from the
# "40% S/LFM+CDN/RME/H:1" line.
$ws.Range("B2").Value = "16% MUR/LWAL+CDN/H:1`n40% S/LFM+CDN/H:1`n15% S+SL/LFM+CDN/H:1`n7% S/LFBR+CDN/H:1`n11% CR/LFM+CDN/H:1`n11% CR+PC/LFM+CDN/H:1"

# Wrap the (now multi-line) text and size the row to fit it.
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 320

# Move the active selection to E1.
$ws.Range("E1").Select() | Out-Null
